# Update column F ("想去人数") values on the "展览" and "全部类型" sheets.
# Both sheets contain the same underlying rows of data, so the same
# row -> new value mapping applies to each sheet.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 6
    4  = 10360
    8  = 7177
    9  = 17
    11 = 203
    13 = 3196
    14 = 36
    15 = 313
    16 = 659
    18 = 1035
    19 = 278
    20 = 67
    21 = 1635
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
